$d = $word.ActiveDocument

# 1) Merge "CREADOR: " + "ALEXANDER BARRIOS" into a single run's text
#    "CREADOR: ALEXANDER BARRIOS" (this also drops the now-unneeded
#    xml:space="preserve" since the replacement text has no leading/
#    trailing space). Find/Replace only rewrites the text of the run it
#    matched, leaving the bookmark and the following run untouched.
$rng = $d.Content
$rng.Find.Execute("CREADOR: ", $true, $false, $false, $false, $false, $true, 1, $false, "CREADOR: ALEXANDER BARRIOS", 2)

# 2) Remove the now-redundant second run that still holds the old
#    "ALEXANDER BARRIOS" text, without disturbing the bookmark that sits
#    right before it.
$creadorPara = $d.Paragraphs(2)
$oldRunStart = $creadorPara.Range.Start + ("CREADOR: ALEXANDER BARRIOS").Length
$oldRunEnd = $creadorPara.Range.End - 1
$d.Range($oldRunStart, $oldRunEnd).Delete()

# 3) Split off a brand-new paragraph right after the CREADOR paragraph and
#    fill it with the "Prueba practica de reigh" text, each misspelled
#    word wrapped in proofErr spell-check markers, bold + en-US like the
#    surrounding heading runs.
$creadorPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(3)

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Prueba</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>práctica</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> de </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>reigh</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$newPara.Range.InsertXML($xmlFrag)
